# Trade #3 closed at 2026-02-17 20:47:46 - unknown UNKNOWN +0.000%
#
# 1) Summary sheet: bump Total Trades / recompute Win Rate %
# 2) All Trades sheet: trade #31 (row 32) transitions OPEN -> CLOSED (early exit),
#    and a brand new trade #64 (row 65) is appended as OPEN.
# 3) MarketMaking sheet: the same new trade #64 is appended as row 32 (this
#    strategy-specific sheet has a different column layout than "All Trades").

$wb = $excel.ActiveWorkbook

# Helper: write a literal date/time-looking string ("YYYY-MM-DD") into a cell
# without Excel's autodetection turning it into a real date serial number.
# Round-tripping it through a text formula + paste-as-values keeps the cell a
# plain string, matching how the rest of this workbook stores its Date column.
function Set-LiteralText {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 31        # Total Trades: 30 -> 31
$summary.Range("B9").Value = 41.94     # Win Rate %: 43.33 -> 41.94

# ---------------------------------------------------------------------------
# 2) All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #31 (row 32) closes out via early exit
$allTrades.Range("G32").Value = 0.16
$allTrades.Range("H32").Value = "CLOSED"
$allTrades.Range("K32").Value = 100.32
$allTrades.Range("L32").Value = "early_exit"
$allTrades.Range("M32").Value = 0.13

# New trade #64 (row 65) opens
$allTrades.Range("A65").Value = 64
Set-LiteralText $allTrades.Range("B65") "2026-02-17"
$allTrades.Range("C65").Value = "20:47:40"
$allTrades.Range("D65").Value = "MarketMaking"
$allTrades.Range("E65").Value = "UP"
$allTrades.Range("F65").Value = 0.16
$allTrades.Range("H65").Value = "OPEN"
$allTrades.Range("I65").Value = 0
$allTrades.Range("J65").Value = 0
$allTrades.Range("K65").Value = 100.32
$allTrades.Range("M65").Value = 0
$allTrades.Range("N65").Value = 0
$allTrades.Range("O65").Value = 0
$allTrades.Range("P65").Value = 0.6
$allTrades.Range("Q65").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# 3) MarketMaking sheet (strategy-specific log; note the different column
#    order vs. "All Trades": L/M/N/O/P/Q here are Entry Slippage, Exit
#    Slippage, Confidence, Entry Reason, Exit Reason, Duration).
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("A32").Value = 64
Set-LiteralText $mm.Range("B32") "2026-02-17"
$mm.Range("C32").Value = "20:47:40"
$mm.Range("D32").Value = "MarketMaking"
$mm.Range("E32").Value = "UP"
$mm.Range("F32").Value = 0.16
$mm.Range("H32").Value = "OPEN"
$mm.Range("I32").Value = 0
$mm.Range("J32").Value = 0
$mm.Range("K32").Value = 100.32
$mm.Range("L32").Value = 0
$mm.Range("M32").Value = 0
$mm.Range("N32").Value = 0.6
$mm.Range("O32").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q32").Value = 0

Write-Output "edit complete"
